# Hortaliza, Feria Lagunitas de Puerto Montt - Ajo
# Insert two new weekly price records (rows 266 and 267), pushing the
# existing rows 266-339 down to 268-341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 266.
$ws.Rows("266:267").Insert()

# New row 266
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44841
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = 100112003
$ws.Range("G266").Value = "Ajo"
$ws.Range("H266").Value = "Chino"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 210
$ws.Range("K266").Value = 20000
$ws.Range("L266").Value = 20000
$ws.Range("M266").Value = 20000
$ws.Range("N266").Value = "$/caja 10 kilos"
$ws.Range("O266").Value = "China"
$ws.Range("P266").Value = 2000
$ws.Range("Q266").Value = 10
$ws.Range("R266").Value = "Hortaliza"

# New row 267
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44841
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = 100112003
$ws.Range("G267").Value = "Ajo"
$ws.Range("H267").Value = "Chino"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 120
$ws.Range("K267").Value = 23000
$ws.Range("L267").Value = 23000
$ws.Range("M267").Value = 23000
$ws.Range("N267").Value = "$/malla 10 kilos"
$ws.Range("O267").Value = "China"
$ws.Range("P267").Value = 2300
$ws.Range("Q267").Value = 10
$ws.Range("R267").Value = "Hortaliza"
